$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename header columns (row 1)
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 2) Title-case the Spanish connector words (de/del/la/las/los/el/y) in state/municipality names
$renames = @(
  @("B7", "Pabellón De Arteaga"),
  @("B8", "Rincón De Romos"),
  @("B9", "San Francisco De Los Romo"),
  @("B10", "San José De Gracia"),
  @("B14", "Playas De Rosarito"),
  @("B36", "Amatenango De La Frontera"),
  @("B39", "Bejucal De Ocampo"),
  @("B41", "Benemérito De Las Américas"),
  @("B47", "Chiapa De Corzo"),
  @("B50", "Comitán De Domínguez"),
  @("B76", "Mazapa De Madero"),
  @("B81", "Ocozocoautla De Espinosa"),
  @("B88", "Salto De Agua"),
  @("B89", "San Cristóbal De Las Casas"),
  @("B123", "Coyame Del Sotol"),
  @("B129", "Guadalupe Y Calvo"),
  @("B131", "Hidalgo Del Parral"),
  @("B143", "San Francisco Del Oro"),
  @("B165", "San Juan De Sabinas"),
  @("B177", "Villa De Álvarez"),
  @("A179", "Ciudad De México"),
  @("B183", "Cuajimalpa De Morelos"),
  @("B210", "Nombre De Dios"),
  @("B214", "Pánuco De Coronado"),
  @("B221", "San Juan De Guadalupe"),
  @("B222", "San Juan Del Río"),
  @("B223", "San Luis Del Cordero"),
  @("B224", "San Pedro Del Gallo"),
  @("A234", "Estado De México"),
  @("B234", "Acambay De Ruíz Castañeda"),
  @("B237", "Almoloya De Alquisiras"),
  @("B238", "Almoloya De Juárez"),
  @("B239", "Almoloya Del Río"),
  @("B245", "Atizapán De Zaragoza"),
  @("B250", "Chapa De Mota"),
  @("B253", "Coacalco De Berriozábal"),
  @("B259", "Ecatepec De Morelos"),
  @("B265", "Ixtapan De La Sal"),
  @("B266", "Ixtapan Del Oro"),
  @("B279", "Naucalpan De Juárez"),
  @("B290", "San Felipe Del Progreso"),
  @("B291", "San Martín De Las Pirámides"),
  @("B292", "San Simón De Guerrero"),
  @("B294", "Soyaniquilpan De Juárez"),
  @("B303", "Tenango Del Aire"),
  @("B304", "Tenango Del Valle"),
  @("B314", "Tlalnepantla De Baz"),
  @("B319", "Valle De Bravo"),
  @("B320", "Valle De Chalco Solidaridad"),
  @("B321", "Villa De Allende"),
  @("B322", "Villa Del Carbón"),
  @("B334", "San Miguel De Allende"),
  @("B335", "Apaseo El Alto"),
  @("B336", "Apaseo El Grande"),
  @("B344", "Dolores Hidalgo Cuna De La Independencia Nacional"),
  @("B348", "Jaral Del Progreso"),
  @("B356", "Purísima Del Rincón"),
  @("B360", "San Diego De La Unión"),
  @("B362", "San Francisco Del Rincón"),
  @("B364", "San Luis De La Paz"),
  @("B366", "Santa Cruz De Juventino Rosas"),
  @("B368", "Silao De La Victoria"),
  @("B373", "Valle De Santiago"),
  @("B379", "Acapulco De Juárez"),
  @("B382", "Ajuchitlán Del Progreso"),
  @("B383", "Alcozauca De Guerrero"),
  @("B386", "Atenango Del Río"),
  @("B387", "Atlamajalcingo Del Monte"),
  @("B389", "Atoyac De Álvarez"),
  @("B390", "Ayutla De Los Libres"),
  @("B393", "Buenavista De Cuéllar"),
  @("B394", "Chilapa De Álvarez"),
  @("B395", "Chilpancingo De Los Bravo"),
  @("B396", "Coahuayutla De José María Izazaga"),
  @("B401", "Coyuca De Benítez"),
  @("B402", "Coyuca De Catalán"),
  @("B406", "Cuetzala Del Progreso"),
  @("B407", "Cutzamala De Pinzón"),
  @("B413", "Huitzuco De Los Figueroa"),
  @("B414", "Iguala De La Independencia"),
  @("B416", "Ixcateopan De Cuauhtémoc"),
  @("B417", "Zihuatanejo De Azueta"),
  @("B419", "La Unión De Isidoro Montes De Oca"),
  @("B422", "Mártir De Cuilapan"),
  @("B435", "Taxco De Alarcón"),
  @("B437", "Técpan De Galeana"),
  @("B439", "Tepecoacuilco De Trujano"),
  @("B441", "Tixtla De Guerrero"),
  @("B444", "Tlapa De Comonfort"),
  @("B456", "Agua Blanca De Iturbide"),
  @("B462", "Atotonilco De Tula"),
  @("B463", "Atotonilco El Grande"),
  @("B469", "Cuautepec De Hinojosa"),
  @("B475", "Huasca De Ocampo"),
  @("B479", "Huejutla De Reyes"),
  @("B482", "Jacala De Ledezma"),
  @("B489", "Mineral Del Chico"),
  @("B490", "Mineral Del Monte"),
  @("B491", "Mixquiahuala De Juárez"),
  @("B492", "Molango De Escamilla"),
  @("B494", "Nopala De Villagrán"),
  @("B495", "Omitlán De Juárez"),
  @("B496", "Pachuca De Soto"),
  @("B499", "Progreso De Obregón"),
  @("B505", "Santiago De Anaya"),
  @("B506", "Santiago Tulantepec De Lugo Guerrero"),
  @("B510", "Tenango De Doria"),
  @("B512", "Tepehuacán De Guerrero"),
  @("B513", "Tepeji Del Río De Ocampo"),
  @("B515", "Tezontepec De Aldama"),
  @("B521", "Tula De Allende"),
  @("B522", "Tulancingo De Bravo"),
  @("B526", "Zacualtipán De Ángeles"),
  @("B527", "Zapotlán De Juárez"),
  @("B532", "Acatlán De Juárez"),
  @("B533", "Ahualulco De Mercado"),
  @("B537", "Atemajac De Brizuela"),
  @("B540", "Atotonilco El Alto"),
  @("B542", "Autlán De Navarro"),
  @("B549", "Concepción De Buenos Aires"),
  @("B555", "Encarnación De Díaz"),
  @("B560", "Huejuquilla El Alto"),
  @("B561", "Ixtlahuacán De Los Membrillos"),
  @("B562", "Ixtlahuacán Del Río"),
  @("B566", "Jilotlán De Los Dolores"),
  @("B570", "La Manzanilla De La Paz"),
  @("B571", "Lagos De Moreno"),
  @("B577", "Ojuelos De Jalisco"),
  @("B582", "San Cristóbal De La Barranca"),
  @("B583", "San Diego De Alejandría"),
  @("B585", "San Juan De Los Lagos"),
  @("B587", "San Martín De Bolaños"),
  @("B589", "San Miguel El Alto"),
  @("B590", "Santa María De Los Ángeles"),
  @("B593", "Talpa De Allende"),
  @("B594", "Tamazula De Gordiano"),
  @("B596", "Techaluta De Montenegro"),
  @("B599", "Teocuitatlán De Corona"),
  @("B600", "Tepatitlán De Morelos"),
  @("B603", "Tizapán El Alto"),
  @("B604", "Tlajomulco De Zúñiga"),
  @("B614", "Unión De San Antonio"),
  @("B615", "Unión De Tula"),
  @("B616", "Valle De Guadalupe"),
  @("B617", "Valle De Juárez"),
  @("B622", "Yahualica De González Gallo"),
  @("B623", "Zacoalco De Torres"),
  @("B626", "Zapotitlán De Vadillo"),
  @("B627", "Zapotlán Del Rey"),
  @("B628", "Zapotlán El Grande"),
  @("B648", "Coalcomán De Vázquez Pallares"),
  @("B650", "Cojumatlán De Régules"),
  @("B713", "Tiquicheo De Nicolás Romero"),
  @("B736", "Coatlán Del Río"),
  @("B743", "Jonacatepec De Leandro Valle"),
  @("B747", "Puente De Ixtla"),
  @("B752", "Tetela Del Volcán"),
  @("B753", "Tlaltizapán De Zapata"),
  @("B765", "Ixtlán Del Río"),
  @("B772", "Santa María Del Oro"),
  @("B792", "Mier Y Noriega"),
  @("B795", "San Nicolás De Los Garza"),
  @("B800", "Acatlán De Pérez Figueroa"),
  @("B805", "Chalcatongo De Hidalgo"),
  @("B806", "Ciénega De Zimatlán"),
  @("B808", "Coicoyán De Las Flores"),
  @("B809", "Constancia Del Rosario"),
  @("B811", "Cuyamecalco Villa De Zaragoza"),
  @("B812", "Fresnillo De Trujano"),
  @("B813", "Guadalupe De Ramírez"),
  @("B814", "Heroica Ciudad De Ejutla De Crespo"),
  @("B815", "Heroica Ciudad De Huajuapan De León"),
  @("B816", "Heroica Ciudad De Tlaxiaco"),
  @("B817", "Huautla De Jiménez"),
  @("B818", "Ixtlán De Juárez"),
  @("B819", "Heroica Ciudad De Juchitán De Zaragoza"),
  @("B825", "Mariscala De Juárez"),
  @("B826", "Mártires De Tacubaya"),
  @("B828", "Miahuatlán De Porfirio Díaz"),
  @("B830", "Nejapa De Madero"),
  @("B831", "Oaxaca De Juárez"),
  @("B832", "Ocotlán De Morelos"),
  @("B833", "Pinotepa De Don Luis"),
  @("B835", "Putla Villa De Guerrero"),
  @("B847", "San Antonino El Alto"),
  @("B863", "San Felipe Jalapa De Díaz"),
  @("B878", "San José Del Progreso"),
  @("B882", "San Juan Bautista Lo De Soto"),
  @("B889", "San Juan Del Río"),
  @("B913", "San Martín De Los Cansecos"),
  @("B916", "San Mateo Del Mar"),
  @("B922", "San Miguel Del Puerto"),
  @("B923", "San Miguel Del Río"),
  @("B925", "San Miguel El Grande"),
  @("B937", "San Pablo Villa De Mitla"),
  @("B940", "San Pedro El Alto"),
  @("B955", "San Pedro Y San Pablo Ayutla"),
  @("B956", "San Pedro Y San Pablo Teposcolula"),
  @("B964", "Santa Ana Del Valle"),
  @("B979", "Santa Inés De Zaragoza"),
  @("B980", "Santa Inés Del Monte"),
  @("B991", "Santa María Jalapa Del Marqués"),
  @("B1027", "Santo Domingo De Morelos"),
  @("B1039", "Tamazulápam Del Espíritu Santo"),
  @("B1040", "Tanetze De Zaragoza"),
  @("B1042", "Tataltepec De Valdés"),
  @("B1043", "Teococuilco De Marcos Pérez"),
  @("B1044", "Teotitlán De Flores Magón"),
  @("B1045", "Teotitlán Del Valle"),
  @("B1047", "Tepelmeme Villa De Morelos"),
  @("B1048", "Heroica Villa Tezoatlán De Segura Y Luna, Cuna De La Independencia De Oaxaca"),
  @("B1049", "Tlacolula De Matamoros"),
  @("B1050", "Tlalixtac De Cabrera"),
  @("B1051", "Totontepec Villa De Morelos"),
  @("B1054", "Villa De Etla"),
  @("B1055", "Villa De Tututepec"),
  @("B1056", "Villa De Zaachila"),
  @("B1058", "Villa Sola De Vega"),
  @("B1060", "Zapotitlán Del Río"),
  @("B1061", "Zimatlán De Álvarez"),
  @("B1077", "Ayotoxco De Guerrero"),
  @("B1081", "Chalchicomula De Sesma"),
  @("B1089", "Chila De La Sal"),
  @("B1096", "Cuapiaxtla De Madero"),
  @("B1100", "Cuayuca De Andrade"),
  @("B1101", "Cuetzalan Del Progreso"),
  @("B1111", "Huehuetlán El Chico"),
  @("B1112", "Huehuetlán El Grande"),
  @("B1117", "Izúcar De Matamoros"),
  @("B1124", "Los Reyes De Juárez"),
  @("B1130", "Palmar De Bravo"),
  @("B1144", "San Nicolás De Los Ranchos"),
  @("B1146", "San Salvador El Seco"),
  @("B1147", "San Salvador El Verde"),
  @("B1153", "Tecali De Herrera"),
  @("B1159", "Tepanco De López"),
  @("B1160", "Tepango De Rodríguez"),
  @("B1161", "Tepatlaxco De Hidalgo"),
  @("B1164", "Tepexi De Rodríguez"),
  @("B1166", "Tetela De Ocampo"),
  @("B1170", "Tlacotepec De Benito Juárez"),
  @("B1194", "Amealco De Bonfil"),
  @("B1196", "Cadereyta De Montes"),
  @("B1202", "Jalpan De Serra"),
  @("B1203", "Landa De Matamoros"),
  @("B1206", "Pinal De Amoles"),
  @("B1209", "San Juan Del Río"),
  @("B1219", "Armadillo De Los Infante"),
  @("B1220", "Axtla De Terrazas"),
  @("B1226", "Ciudad Del Maíz"),
  @("B1236", "Mexquitic De Carmona"),
  @("B1241", "San Ciro De Acosta"),
  @("B1247", "Santa María Del Río"),
  @("B1249", "Soledad De Graciano Sánchez"),
  @("B1257", "Tanquián De Escobedo"),
  @("B1261", "Villa De Arista"),
  @("B1262", "Villa De Arriaga"),
  @("B1263", "Villa De Guadalupe"),
  @("B1264", "Villa De La Paz"),
  @("B1265", "Villa De Ramos"),
  @("B1266", "Villa De Reyes"),
  @("B1310", "Jalpa De Méndez"),
  @("B1341", "Soto La Marina"),
  @("B1354", "Ixtacuixtla De Mariano Matamoros"),
  @("B1357", "Nanacamilpa De Mariano Arista"),
  @("B1360", "Papalotla De Xicohténcatl"),
  @("B1362", "San Pablo Del Monte"),
  @("B1363", "Sanctórum De Lázaro Cárdenas"),
  @("B1370", "Tetla De La Solidaridad"),
  @("B1385", "Alto Lucero De Gutiérrez Barrios"),
  @("B1389", "Amatlán De Los Reyes"),
  @("B1400", "Boca Del Río"),
  @("B1402", "Camarón De Tejeda"),
  @("B1405", "Castillo De Teayo"),
  @("B1407", "Cazones De Herrera"),
  @("B1425", "Cosamaloapan De Carpio"),
  @("B1426", "Cosautlán De Carvajal"),
  @("B1442", "Hueyapan De Ocampo"),
  @("B1443", "Ignacio De La Llave"),
  @("B1447", "Ixhuacán De Los Reyes"),
  @("B1448", "Ixhuatlán De Madero"),
  @("B1449", "Ixhuatlán Del Café"),
  @("B1450", "Ixhuatlán Del Sureste"),
  @("B1460", "Juchique De Ferrer"),
  @("B1464", "Las Vigas De Ramírez"),
  @("B1465", "Lerdo De Tejada"),
  @("B1468", "Martínez De La Torre"),
  @("B1470", "Medellín De Bravo"),
  @("B1474", "Nanchital De Lázaro Cárdenas Del Río"),
  @("B1482", "Ozuluama De Mascareñas"),
  @("B1486", "Paso De Ovejas"),
  @("B1487", "Paso Del Macho"),
  @("B1491", "Poza Rica De Hidalgo"),
  @("B1499", "Sayula De Alemán"),
  @("B1503", "Soledad De Doblado"),
  @("B1509", "Tatahuicapan De Juárez"),
  @("B1535", "Vega De Alatorre"),
  @("B1543", "Zontecomatlán De López Y Fuentes"),
  @("B1558", "Concepción Del Oro"),
  @("B1568", "Jiménez Del Teul"),
  @("B1574", "Mezquital Del Oro"),
  @("B1578", "Noria De Ángeles"),
  @("B1587", "Teúl De González Ortega"),
  @("B1588", "Tlaltenango De Sánchez Román"),
  @("B1591", "Villa De Cos")
)

foreach ($pair in $renames) {
    $ws.Range($pair[0]).Value = $pair[1]
}

# 3) Remove the trailing blank row and footer/metadata rows (old rows 1599-1604)
$ws.Rows("1599:1604").Delete()

